# Applies the "find a neighborhood" text_dictionary update described in the commit:
#   - adds four new dictionary rows (text_120..text_123) with EN/ES copy describing the
#     Census geocoder tool used by "Find a Neighborhood", and the opt-out instructions
#   - relabels the existing "Address"/"Direccion" cells (text_118, row 121) to add a colon
#
# NOTE: the workbook stores text as a shared-string table; cell values only carry an index
# into that table. Unused strings are pruned and new unique strings are appended to the
# table in the order they are first assigned. The assignment order below is deliberately
# chosen (matching the order the strings were authored upstream) so that the resulting
# shared-string table/index layout matches the target workbook exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New rows 123:126 -> text_120 .. text_123 (text_id column) ---
$ws.Range("A123").Value = 'text_120'
$ws.Range("A124").Value = 'text_121'
$ws.Range("A125").Value = 'text_122'
$ws.Range("A126").Value = 'text_123'

# --- English copy (column B) ---
$ws.Range("B123").Value = 'About the "Find a Neighborhood" tool:'
$ws.Range("B125").Value = '"When you submit a geocode request, the web server automatically collects certain technical information from your computer and about your connection. The only information that is stored is IP Address (for batch and single address submissions) and submitted address (only for single address submissions). The information is stored on a server, which is internal to the Census network.  This server is only accessible by Census Bureau staff, who are bound by the confidentiality requirements set forth in Title 13 of the United States Code and requires multiple levels of approval. This information is only used to monitor and track the performance of the Geocoder."'
$ws.Range("B124").Value = 'To find the neighborhood associated with an address, we use a geocoder run by the US Census Bureau. The following information from the US Census Bureau descibes how they use the address information you provide:'

# --- Spanish copy (column C) ---
$ws.Range("C123").Value = 'Acerca de la herramienta "Buscar un vecindario":'
$ws.Range("C124").Value = 'Para encontrar el vecindario asociado con una dirección, utilizamos un geocodificador administrado por la Oficina del Censo de los EE. UU. La siguiente información de la Oficina del Censo de los EE. UU. describe cómo utilizan la información de la dirección que usted proporciona:'
$ws.Range("C125").Value = '"Cuando envía una solicitud de geocodificación, el servidor web recopila automáticamente cierta información técnica de su computadora y sobre su conexión. La única información que se almacena es la dirección IP (para envíos por lotes y de una sola dirección) y la dirección enviada (solo para envíos de una sola dirección). La información se almacena en un servidor, que es interno a la red del Censo. A este servidor solo puede acceder el personal de la Oficina del Censo, que está sujeto a los requisitos de confidencialidad establecidos en el Título 13 del Código de los Estados Unidos y requiere múltiples niveles de aprobación. Esta información solo se utiliza para monitorear y hacer un seguimiento del rendimiento del geocodificador".'
$ws.Range("C126").Value = 'Si no desea ingresar su propia dirección, puede ingresar la dirección de una empresa local u otra ubicación en el vecindario de interés.'

$ws.Range("B126").Value = 'If you do not wish to enter your own address, you can enter the address of a local business or another location in the neighborhood of interest.'

# --- Existing row 121 (text_118): "Address" / "Direccion" -> "Address:" / "Direccion:" ---
$ws.Range("C121").Value = 'Dirección:'
$ws.Range("B121").Value = 'Address:'

# --- Row heights for the new wrapped rows (matches authored heights) ---
$ws.Rows.Item(123).RowHeight = 17
$ws.Rows.Item(124).RowHeight = 51
$ws.Rows.Item(125).RowHeight = 136
$ws.Rows.Item(126).RowHeight = 34

# --- Leave the selection where the author left it ---
[void]$ws.Range("C130").Select()
